$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the commit diff.
# Columns D/E hold numeric-looking / percentage text; prefix with a literal
# apostrophe so Excel stores them as text (quote-prefixed), matching the
# original inlineStr cell type rather than being coerced to a number.

$ws.Range("D2").Value = '''26.475.55'
$ws.Range("E2").Value = '''  -3.32%  '
$ws.Range("D3").Value = '''1.668.69'
$ws.Range("E3").Value = '''  -2.47%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '''  +0.12%  '
$ws.Range("D5").Value = '''218.83'
$ws.Range("E5").Value = '''  -2.36%  '
$ws.Range("D6").Value = '''0.5161'
$ws.Range("E6").Value = '''  -3.06%  '
$ws.Range("E7").Value = '''  +0.23%  '
$ws.Range("D8").Value = '''0.06474'
$ws.Range("E8").Value = '''  -1.81%  '
$ws.Range("E9").Value = '''  -3.10%  '
$ws.Range("D10").Value = '''20.01'
$ws.Range("E10").Value = '''  -4.18%  '
$ws.Range("D11").Value = '''0.07666'
$ws.Range("E11").Value = '''  +0.38%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '''1.677.97'
$ws.Range("E12").Value = '''  -2.23%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.342'
$ws.Range("E13").Value = '''  -5.09%  '
$ws.Range("D14").Value = '''1.898.96'
$ws.Range("E14").Value = '''  -2.31%  '
$ws.Range("D15").Value = '''0.5547'
$ws.Range("E15").Value = '''  -3.17%  '
$ws.Range("D16").Value = '''0.0₅8059'
$ws.Range("E16").Value = '''  -1.42%  '
$ws.Range("D17").Value = '''64.73'
$ws.Range("E17").Value = '''  -4.64%  '
$ws.Range("D18").Value = '''26.500.59'
$ws.Range("E18").Value = '''  -3.15%  '
$ws.Range("D19").Value = '''1.006'
$ws.Range("E19").Value = '''  +0.27%  '
$ws.Range("D20").Value = '''210.84'
$ws.Range("E20").Value = '''  -2.35%  '
$ws.Range("D21").Value = '''4.423'
$ws.Range("E21").Value = '''  -5.36%  '
$ws.Range("E22").Value = '''  -2.94%  '
$ws.Range("D23").Value = '''5.894'
$ws.Range("E23").Value = '''  -1.36%  '
$ws.Range("E24").Value = '''  +0.28%  '
$ws.Range("D25").Value = '''144.94'
$ws.Range("E25").Value = '''  +2.32%  '
$ws.Range("E26").Value = '''  -1.58%  '
$ws.Range("D27").Value = '''0.1163'
$ws.Range("E27").Value = '''  -4.52%  '
$ws.Range("D28").Value = '''7.021'
$ws.Range("E28").Value = '''  -3.48%  '
$ws.Range("D29").Value = '''15.79'
$ws.Range("E29").Value = '''  -3.17%  '
$ws.Range("D30").Value = '''0.05233'
$ws.Range("E30").Value = '''  -3.37%  '
$ws.Range("D31").Value = '''1.263'
$ws.Range("E31").Value = '''  -2.42%  '
$ws.Range("D32").Value = '''3.371'
$ws.Range("E32").Value = '''  -3.92%  '
$ws.Range("D33").Value = '''3.228'
$ws.Range("E33").Value = '''  -5.92%  '
$ws.Range("D34").Value = '''1.578'
$ws.Range("E34").Value = '''  -4.10%  '
$ws.Range("D35").Value = '''2.760'
$ws.Range("E35").Value = '''  -4.18%  '
$ws.Range("D36").Value = '''2.375'
$ws.Range("E36").Value = '''  -1.84%  '
$ws.Range("D37").Value = '''0.9250'
$ws.Range("E37").Value = '''  -2.57%  '
$ws.Range("D38").Value = '''0.5738'
$ws.Range("E38").Value = '''  -2.08%  '
$ws.Range("D39").Value = '''1.164.38'
$ws.Range("E39").Value = '''  +11.14%  '
$ws.Range("E40").Value = '''  -1.99%  '
$ws.Range("D41").Value = '''1.006'
$ws.Range("E41").Value = '''  +0.22%  '
$ws.Range("D42").Value = '''0.8450'
$ws.Range("E42").Value = '''  +0.43%  '
$ws.Range("E43").Value = '''  -3.85%  '
$ws.Range("D44").Value = '''100.28'
$ws.Range("E44").Value = '''  -0.56%  '
$ws.Range("D45").Value = '''1.807.56'
$ws.Range("E45").Value = '''  -2.35%  '
$ws.Range("D46").Value = '''0.0₈113'
$ws.Range("E46").Value = '''  -2.51%  '
$ws.Range("D47").Value = '''0.4497'
$ws.Range("E47").Value = '''  -0.24%  '
$ws.Range("D48").Value = '''56.09'
$ws.Range("E48").Value = '''  -3.26%  '
$ws.Range("E49").Value = '''  +0.20%  '
$ws.Range("D50").Value = '''7.972'
$ws.Range("E50").Value = '''  -1.59%  '
$ws.Range("D51").Value = '''0.05120'
$ws.Range("E51").Value = '''  -2.39%  '
